# "Generate Report for Handoff"
# A new handoff package was generated (new content guid + new content
# hash), so every reference to the old guid/hash is refreshed, the
# handoff timestamps advance, and the now-stale "handback" columns are
# reset because the freshly generated handoff has not been handed back
# yet.

$newGuid = "8d264554-8f2f-4fdc-991b-1d8b7cc2fb08"
$newHash = "3ad874ef0452ccdc69a841c4b3e85f13fe3a9283"

$newGenerateDate = "2016-09-06 21:21:00"
$newHandoffDateZh = "2016-09-06 21:20:55"
$resetHandbackDate = "0001-01-01 00:00:00"

# The hyperlink *targets* are historical GitHub blob URLs that keep
# pointing at the commit that was current when they were first created -
# they are left untouched; only the cell text / hyperlink display text
# moves on to the new guid.
$overviewAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d76fb1add516af8839da5f2145d6e4ecc92013b6/e2e/e414559a-85d2-4c60-8b29-5c9aa639a168.md"
$zhcnAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d76fb1add516af8839da5f2145d6e4ecc92013b6/e2e/e414559a-85d2-4c60-8b29-5c9aa639a168.md"
$dedeAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d76fb1add516af8839da5f2145d6e4ecc92013b6/e2e/e414559a-85d2-4c60-8b29-5c9aa639a168.md"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("G2").Value = $newGenerateDate

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), $overviewAddr, $null, $null, "e2e\$newGuid.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$ws.Range("H2").Value = $newHandoffDateZh

# Freshly (re)generated handoff -> nothing has been targeted/handed back
# against it yet, so these clear out (and lose their old hyperlink).
$ws.Range("I2").Value = ""
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = $resetHandbackDate

$ws.Hyperlinks.Add($ws.Range("A2"), $zhcnAddr, $null, $null, "$newGuid.md")

$ws.Columns.Item(9).ColumnWidth = 18.6506053379604 - (5/6)
$ws.Columns.Item(10).ColumnWidth = 21.7054770333426 - (5/6)

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$ws.Range("H2").Value = $newGenerateDate

$ws.Range("I2").Value = ""
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = $resetHandbackDate

$ws.Hyperlinks.Add($ws.Range("A2"), $dedeAddr, $null, $null, "$newGuid.md")

$ws.Columns.Item(9).ColumnWidth = 18.6506053379604 - (5/6)
$ws.Columns.Item(10).ColumnWidth = 21.7054770333426 - (5/6)
